{"js": "// Relatorio-Modulo2.docx \u2014 \"pequena correcao no relat\u00f3rio\"\n//\n// 1) Fix the \"Singlethon\" -> \"Singleton\" typo (occurs twice).\n// 2) Fix \"diminuimos\" -> \"diminu\u00edmos\" (missing accent).\n// 3) Fix \"manipular\" -> \"manipulador\" (\"criar o manipular de reposit\u00f3rio\"\n//    should read \"criar o manipulador de reposit\u00f3rio\").\n// 4) Fix the missing space after the period in \"usu\u00e1rio.H\u00e1\" -> \"usu\u00e1rio. H\u00e1\".\n// 5) Remove the (both/justify) alignment on the \"J\u00e1 o Factory...\" paragraph\n//    so it falls back to the default (left) alignment.\n// 6) Add a primary/default header containing the \"#interna\" sensitivity\n//    label text (this also provisions the even/first header & footer\n//    parts and footnotes/endnotes referenced from the section).\n\nconst body = context.document.body;\n\n// --- 1) Singlethon -> Singleton (both occurrences) ---------------------\nconst singlethon = body.search(\"Singlethon\", { matchCase: true });\nsinglethon.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < singlethon.items.length; i++) {\n  singlethon.items[i].insertText(\"Singleton\", Word.InsertLocation.replace);\n}\n\n// --- 2) diminuimos -> diminu\u00edmos ----------------------------------------\nconst diminuimos = body.search(\"diminuimos\", { matchCase: true });\ndiminuimos.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < diminuimos.items.length; i++) {\n  diminuimos.items[i].insertText(\"diminu\u00edmos\", Word.InsertLocation.replace);\n}\n\n// --- 3) manipular -> manipulador (only the repository-handler mention) -\nconst manipular = body.search(\"criar o manipular de reposit\u00f3rio\", { matchCase: true });\nmanipular.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < manipular.items.length; i++) {\n  manipular.items[i].insertText(\"criar o manipulador de reposit\u00f3rio\", Word.InsertLocation.replace);\n}\n\n// --- 4) usu\u00e1rio.H\u00e1 -> usu\u00e1rio. H\u00e1 ---------------------------------------\n// Search a slightly wider span (instead of just \"usu\u00e1rio.H\u00e1\") so the\n// replacement also swallows the now-stale spelling/grammar proofing marks\n// that flagged the missing space \u2014 once the space is added the sentence is\n// no longer flagged, so those <w:proofErr/> marks should disappear too.\nconst missingSpace = body.search(\"entrada do usu\u00e1rio.H\u00e1 a utiliza\u00e7\u00e3o\", { matchCase: true });\nmissingSpace.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < missingSpace.items.length; i++) {\n  missingSpace.items[i].insertText(\"entrada do usu\u00e1rio. H\u00e1 a utiliza\u00e7\u00e3o\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// --- 5) Drop the justified alignment on the \"J\u00e1 o Factory...\" paragraph\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"J\u00e1 o Factory\") !== -1) {\n    paragraphs.items[i].alignment = Word.Alignment.left;\n    break;\n  }\n}\nawait context.sync();\n\n// --- 6) Add the \"#interna\" sensitivity-label header ---------------------\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst primaryHeader = sections.items[0].getHeader(Word.HeaderFooterType.primary);\nprimaryHeader.insertText(\"#interna\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Relatorio-Modulo2.docx \u2014 \"pequena correcao no relat\u00f3rio\"\n#\n# 1) Fix the \"Singlethon\" -> \"Singleton\" typo (occurs twice).\n# 2) Fix \"diminuimos\" -> \"diminu\u00edmos\" (missing accent).\n# 3) Fix \"manipular\" -> \"manipulador\" (\"criar o manipular de reposit\u00f3rio\"\n#    should read \"criar o manipulador de reposit\u00f3rio\").\n# 4) Fix the missing space after the period in \"usu\u00e1rio.H\u00e1\" -> \"usu\u00e1rio. H\u00e1\"\n#    (a wider find/replace span is used so the stale spelling/grammar\n#    <w:proofErr/> marks around the old typo go away with it).\n# 5) Remove the (both/justify) alignment on the \"J\u00e1 o Factory...\" paragraph\n#    so it falls back to the default (left) alignment.\n# 6) Add a primary/default header containing the \"#interna\" sensitivity\n#    label text (this also provisions the even/first header & footer\n#    parts and footnotes/endnotes referenced from the section).\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# --- 1) Singlethon -> Singleton (both occurrences) ----------------------\nReplace-AllText \"Singlethon\" \"Singleton\"\n\n# --- 2) diminuimos -> diminu\u00edmos -----------------------------------------\nReplace-AllText \"diminuimos\" \"diminu\u00edmos\"\n\n# --- 3) manipular -> manipulador (only the repository-handler mention) --\nReplace-AllText \"criar o manipular de reposit\u00f3rio\" \"criar o manipulador de reposit\u00f3rio\"\n\n# --- 4) usu\u00e1rio.H\u00e1 -> usu\u00e1rio. H\u00e1 ----------------------------------------\nReplace-AllText \"entrada do usu\u00e1rio.H\u00e1 a utiliza\u00e7\u00e3o\" \"entrada do usu\u00e1rio. H\u00e1 a utiliza\u00e7\u00e3o\"\n\n# --- 5) Drop the justified alignment on the \"J\u00e1 o Factory...\" paragraph -\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*J\u00e1 o Factory*\") {\n        $p.Alignment = 0   # wdAlignParagraphLeft\n        break\n    }\n}\n\n# --- 6) Add the \"#interna\" sensitivity-label header ----------------------\n$sec = $d.Sections.Item(1)\n$hdr = $sec.Headers.Item(1)   # wdHeaderFooterPrimary\n$hdr.Range.Text = \"#interna\"\n\nWrite-Output \"done\"\n"}
